$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 57
$ws.Range("I2").Value = 164
$ws.Range("J2").Value = 580
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 179
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 100
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 58
$ws.Range("T2").Value = 101
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 958
$ws.Range("X2").Value = 886
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 6
